# Update the "相談件数" (consultation count) sheet with one new day's data row.
#
# Before:  sheetData ends with row 101 = the footnote row ("※4/8...")
# After:   a new data row is inserted as row 101 (date 2020-05-05 / 43956,
#          with values 321 / 33785 / 0 / 6958) and the footnote row is
#          pushed down to row 102. The worksheet dimension, the sheet's
#          local Print_Area defined name, and the frozen-pane selection
#          are updated to reflect the extra row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")

# Insert a new row just above the old footnote row (currently row 101),
# pushing the footnote down to row 102. The new row inherits the number
# formatting/styles from the row above it (row 100), matching the target.
$ws.Rows.Item(101).Insert()

# Populate the newly inserted row 101 with the new day's figures.
$ws.Range("A101").Value = 43956
$ws.Range("B101").Value = 321
$ws.Range("C101").Value = 33785
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 6958

# Update the sheet-local "_xlnm.Print_Area" defined name so the print
# area grows by one row (was $A$1:$E$102, now $A$1:$E$103).
$printArea = $wb.Names.Item("相談件数!Print_Area")
$printArea.RefersTo = "=相談件数!`$A`$1:`$E`$103"

# Move the selection in the frozen (bottom-right) pane to the cell below
# the footnote label, matching where the active cell ends up after the
# row insert/edit (B102).
$null = $ws.Range("B102").Select()
